$wb = $excel.ActiveWorkbook

# --- Sheet "Create": add a second column (ErrorMsg1 / Enter Role Name) ---
# and rename the sample role from "Palakadmin" to "PalakAdmin"
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Activate()
$wsCreate.Range("B1").Value = "ErrorMsg1"
$wsCreate.Range("A2").Value = "PalakAdmin"
$wsCreate.Range("B2").Value = "Enter Role Name"
$wsCreate.Range("B1").Select()

# --- Sheet "Edit": update the sample role / updated name / reason ---
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Activate()
# Leading apostrophe keeps the cell's existing "quote prefix" text style
$wsEdit.Range("A2").Value = "'PalakAdmin"
$wsEdit.Range("B2").Value = "PalakAdmin Updated"
$wsEdit.Range("C2").Value = "Modified"
$wsEdit.Range("A1:C2").Select()

# --- Sheet "Delete": update the sample role / reason, make this the active tab ---
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Activate()
$wsDelete.Range("B2").Value = "deleted"
$wsDelete.Range("A2").Value = "Empire"
$wsDelete.Range("A2").Select()
